# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
# (commit: "Updated symbol list on Fri Jan 20 05:33:26 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'291.82"
$ws.Range("E2").Value = "'0.37%"
$ws.Range("E3").Value = "'0.52%"
$ws.Range("D4").Value = "'4.921"
$ws.Range("E4").Value = "'-0.20%"
$ws.Range("E5").Value = "'2.80%"
$ws.Range("D6").Value = "'2.207"
$ws.Range("E6").Value = "'22.90%"
$ws.Range("E7").Value = "'0.84%"
$ws.Range("D8").Value = "'3.754"
$ws.Range("E8").Value = "'0.17%"
$ws.Range("D9").Value = "'0.9122"
$ws.Range("E9").Value = "'1.85%"
$ws.Range("D10").Value = "'0.08922"
$ws.Range("E10").Value = "'15.72%"
$ws.Range("D11").Value = "'0.1688"
$ws.Range("E11").Value = "'2.36%"
$ws.Range("D12").Value = "'0.08265"
$ws.Range("E12").Value = "'2.93%"
$ws.Range("D13").Value = "'0.03109"
$ws.Range("E13").Value = "'1.46%"
$ws.Range("D14").Value = "'0.09991"
$ws.Range("E14").Value = "'-0.23%"
$ws.Range("D15").Value = "'0.001498"
$ws.Range("E15").Value = "'-0.12%"
$ws.Range("D16").Value = "'0.005852"
$ws.Range("E16").Value = "'2.68%"
$ws.Range("D17").Value = "'3.497"
$ws.Range("E17").Value = "'0.60%"
$ws.Range("D18").Value = "'2.087"
$ws.Range("E18").Value = "'0.33%"
$ws.Range("D20").Value = "'0.1298"
$ws.Range("E20").Value = "'1.93%"
$ws.Range("D21").Value = "'3.983"
$ws.Range("E21").Value = "'-1.60%"
$ws.Range("D22").Value = "'0.2190"
$ws.Range("E22").Value = "'9.57%"
$ws.Range("E23").Value = "'0.98%"
$ws.Range("E24").Value = "'0.18%"
$ws.Range("D25").Value = "'0.004580"
$ws.Range("E25").Value = "'14.31%"
$ws.Range("D26").Value = "'0.0001302"
$ws.Range("E26").Value = "'4.16%"
$ws.Range("D27").Value = "'0.0003402"
$ws.Range("D39").Value = "'0.01592"
$ws.Range("E39").Value = "'-0.59%"
$ws.Range("D40").Value = "'0.04465"
$ws.Range("E40").Value = "'1.77%"
$ws.Range("D41").Value = "'0.007367"
$ws.Range("E41").Value = "'-0.03%"
$ws.Range("D42").Value = "'0.009567"
$ws.Range("E42").Value = "'24.70%"
$ws.Range("E43").Value = "'1.49%"
$ws.Range("D44").Value = "'0.002334"
$ws.Range("E44").Value = "'12.73%"
$ws.Range("D45").Value = "'0.009144"
$ws.Range("E45").Value = "'-0.70%"
$ws.Range("D46").Value = "'0.00006106"
$ws.Range("E46").Value = "'3.14%"
$ws.Range("E47").Value = "'0.15%"
$ws.Range("D48").Value = "'2.113"
$ws.Range("E48").Value = "'-5.96%"
$ws.Range("D49").Value = "'0.002005"
$ws.Range("E49").Value = "'-33.20%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.15%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.15%"
